$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for columns B..Q, identical across all data rows (2..26)
$newValues = @(
    0.5692387035740197,
    0.2355707868591003,
    0.4408509555957502,
    0.3174388163398759,
    0.4793110825313374,
    0.2557185134775685,
    0.4537982025430372,
    0.4118254736402926,
    0.2354585389209425,
    0.3236420062806176,
    0.2823828768330753,
    0.5056861808251918,
    0.06015717143422494,
    0.5272142884576241,
    28.72735599237364,
    44.57274171566024
)

for ($row = 2; $row -le 26; $row++) {
    for ($col = 2; $col -le 17; $col++) {
        $ws.Cells.Item($row, $col).Value = $newValues[$col - 2]
    }
}
